# QuestTable.xlsx edit: update quest goal-argument data (columns C/D) for
# rows 2-4 on Sheet1, and move the active selection to D4 (matching the
# author's working selection at save time).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 2 (quest index 1): D2 goalArgument 1 -> 2
$ws.Cells.Item(2, 4).Value = 2

# Row 3 (quest index 2): C3 goalType 2 -> 3, D3 goalArgument 2 -> 4
$ws.Cells.Item(3, 3).Value = 3
$ws.Cells.Item(3, 4).Value = 4

# Row 4 (quest index 3): C4 goalType 3 -> 5, D4 goalArgument 3 -> 6
$ws.Cells.Item(4, 3).Value = 5
$ws.Cells.Item(4, 4).Value = 6

# Reposition the visible selection to D4 (matches the saved view state)
$ws.Range("D4").Select()
